$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update roster: "Elin Min" steps down -> "TBD" takes her Senior VP slot,
# and new member "Julia Lin" fills the previously-TBD Pledge Education VP slot.
$ws.Range("A3").Value = "TBD"
$ws.Range("A4").Value = "Julia Lin"

# Move active selection to A4
$ws.Range("A4").Select()
